# Update the "Platform Coverage" sheet with revised coverage figures and
# bring the newly-edited columns into view / selection, matching the
# author's on-screen state when the workbook was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Platform Coverage")
$ws.Activate()

# Row 2 (pre-SAC): max age bumped 14 -> 15, and the 2026 coverage (P2)
# raised from 0.6 to 0.8.
$ws.Range("G2").Value = 15
$ws.Range("P2").Value = 0.8

# Row 3 (SAC): max age bumped 49 -> 50, and a 2026 coverage value (P3) of
# 0.5 is now supplied (previously blank).
$ws.Range("G3").Value = 50
$ws.Range("P3").Value = 0.5

# Row 4 (adults): 2026 coverage value (P4) of 0.5 is now supplied
# (previously blank); G4 (65) is unchanged.
$ws.Range("P4").Value = 0.5

# Reflect the author's final on-screen view: scrolled right so column G
# is the first visible column, with P2:P4 selected (P2 active).
$excel.ActiveWindow.ScrollColumn = 7
$ws.Range("P2:P4").Select()
